$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 8, shifting existing rows
# 8-61 down to 9-62 (old row 61 becomes row 62).
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly price entry.
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Femacal de La Calera"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44764
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 100112035
$ws.Range("G8").Value = "Bruselas (repollito)"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 45
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("N8").Value = "$/malla 15 kilos"
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 1000
$ws.Range("Q8").Value = 15
$ws.Range("R8").Value = "Hortaliza"
